$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update B26: phone number value
$ws.Range("B26").Value = "9036544535"

# Update C26: password value
$ws.Range("C26").Value = "akhi2506"

# Update selected cell in the sheet view
$ws.Range("E23").Select()
